# PlantillaLineasFacturaRecibida.xlsx
# Rename two header columns in the "Table1" table on sheet "PaginaDetalle":
#   A1: IdTipoDocumento -> TipoDocumento
#   T1: CodigoImpuesto  -> CodigoEtax
# (Supports new import flows: import by XML and by Email - facturas@etaxcr.com)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the table header cells directly - this renames the ListObject's
# columns (xl/tables/table1.xml) and updates the shared strings / cell
# values used by the header row in one go.
$ws.Range("A1").Value = "TipoDocumento"
$ws.Range("T1").Value = "CodigoEtax"

# Reflect the new selection/active cell left in the sheet when it was saved.
$ws.Range("S7").Select()
